$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting existing rows 82-184 down to 83-185
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with the new record's data
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value = "Los Lagos"
$ws.Range("D82").Value = 44546
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 100112021
$ws.Range("G82").Value = "Ají"
$ws.Range("H82").Value = "Inferno"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 70
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 20000
$ws.Range("M82").Value = 20000
$ws.Range("N82").Value = "$/caja 12 kilos"
$ws.Range("O82").Value = "Región de Arica y Parinacota"
$ws.Range("P82").Value = 1667
$ws.Range("Q82").Value = 12
$ws.Range("R82").Value = "Hortaliza"
